# Apply cell-value updates to the crypto price/volume sheet
# (symbol-list refresh, per the commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E columns hold numeric-looking text (prices, % changes) that must
# stay text, so force the Text number format before writing each value.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "326.54"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.38%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.10"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.567"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.16%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08028"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-4.10%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.295"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.14%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.888"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.41%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9435"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.48%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1147"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-7.54%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1832"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.26%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09691"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.52%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04366"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.32%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.32%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001271"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.23%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04220"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.14%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005999"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.63%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.607"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.49%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.19%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.590"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.15%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.05%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2655"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.52%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001251"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.45%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004490"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.38%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.04%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003996"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.06%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02608"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-7.19%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05420"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-8.65%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-4.04%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1395"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.06%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007278"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-19.50%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002020"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.92%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008829"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-14.77%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006928"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.30%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.04%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003632"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "13.39%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002274"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.06%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.04%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.04%"
